$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("118÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "919÷6=", 2) | Out-Null
$d.Content.Find.Execute("479÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "585÷9=", 2) | Out-Null
$d.Content.Find.Execute("147÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "582÷2=", 2) | Out-Null
$d.Content.Find.Execute("490÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "372÷4=", 2) | Out-Null
$d.Content.Find.Execute("953÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "843÷4=", 2) | Out-Null
$d.Content.Find.Execute("271÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "983÷4=", 2) | Out-Null
$d.Content.Find.Execute("153÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "964÷6=", 2) | Out-Null
$d.Content.Find.Execute("200÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "583÷3=", 2) | Out-Null
$d.Content.Find.Execute("402÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "289÷2=", 2) | Out-Null
$d.Content.Find.Execute("951÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "959÷3=", 2) | Out-Null
$d.Content.Find.Execute("893÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "918÷8=", 2) | Out-Null
$d.Content.Find.Execute("153÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "853÷7=", 2) | Out-Null
$d.Content.Find.Execute("170÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "763÷2=", 2) | Out-Null
$d.Content.Find.Execute("757÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "611÷4=", 2) | Out-Null
$d.Content.Find.Execute("493÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "387÷4=", 2) | Out-Null
$d.Content.Find.Execute("650÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "371÷9=", 2) | Out-Null
$d.Content.Find.Execute("642÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "701÷6=", 2) | Out-Null
$d.Content.Find.Execute("187÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "148÷3=", 2) | Out-Null
$d.Content.Find.Execute("103÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷3=", 2) | Out-Null
$d.Content.Find.Execute("783÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "444÷4=", 2) | Out-Null
$d.Content.Find.Execute("467÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "359÷5=", 2) | Out-Null
$d.Content.Find.Execute("403÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "639÷4=", 2) | Out-Null
$d.Content.Find.Execute("300÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "262÷9=", 2) | Out-Null
$d.Content.Find.Execute("226÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷3=", 2) | Out-Null
$d.Content.Find.Execute("445÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "768÷4=", 2) | Out-Null
